# issue #5: add legislator_id, name, date into dataframe
# The "股票" (stocks) sheet is the 4th worksheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1): new columns H/I/J ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Copy the header style (bold, centered, bordered) from the existing G1
# header cell onto the three new header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-7): date / legislator_name / legislator_id ---
for ($r = 2; $r -le 7; $r++) {
    $dateCell = $ws.Cells.Item($r, 8)

    # Force text format before assigning, otherwise Excel auto-converts a
    # "YYYY-MM-DD" literal into a date serial number instead of keeping it
    # as the literal string "2011-11-23".
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2011-11-23"

    $ws.Cells.Item($r, 9).Value = "陳明文"
    $ws.Cells.Item($r, 10).Value = 828
}

# Re-apply the plain data-row style (matching column G) to the date column
# so it doesn't keep the custom "@" text format - value is already a
# literal string at this point so it won't get re-parsed as a date.
$ws.Range("G2").Copy()
$ws.Range("H2:H7").PasteSpecial(-4122)
